# Updates the "cryptos" price/volume snapshot (GitHub Actions refresh).
# Column D ("Price") and E ("Volume(1h)") are plain text cells (not numbers),
# so for values that look numeric we briefly force a Text number format
# before assigning, then restore the default "Normal" style so the cell's
# formatting is unchanged from before (only its stored text differs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.793.07'

$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '1.630.44'

$ws.Range("E3").Value = '  -0.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  -0.68%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.60'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.501'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -0.66%  '

$ws.Range("E8").Value = '  -0.55%  '

$ws.Range("E9").Value = '  -0.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.58'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("E11").Value = '  +1.05%  '

$ws.Range("B12").Value = 'Polkadot'

$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.26'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.51%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'

$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'

$ws.Range("D13").Value = '1.855.89'

$ws.Range("E13").Value = '  -0.06%  '

$ws.Range("B14").Value = 'WrappedEther'

$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'

$ws.Range("D14").Value = '1.613.48'

$ws.Range("E14").Value = '  -1.07%  '

$ws.Range("E15").Value = '  +0.56%  '

$ws.Range("D16").Value = '0.0₃0762'

$ws.Range("E16").Value = '  -0.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.96'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -0.40%  '

$ws.Range("D18").Value = '25.780.47'

$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.999'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -0.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.45'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +0.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.81'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -0.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.93'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("E23").Value = '  +1.02%  '

$ws.Range("E24").Value = '  +1.98%  '

$ws.Range("E25").Value = '  -0.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.73'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +1.74%  '

$ws.Range("E27").Value = '  +3.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.85'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +0.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.50'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -0.19%  '

$ws.Range("E30").Value = '  -0.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0495'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +1.63%  '

$ws.Range("E32").Value = '  -0.03%  '

$ws.Range("E33").Value = '  -0.24%  '

$ws.Range("E34").Value = '  +0.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.909'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +1.37%  '

$ws.Range("D37").Value = '1.140.17'

$ws.Range("E37").Value = '  +2.99%  '

$ws.Range("E38").Value = '  -2.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.543'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -0.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0156'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -0.46%  '

$ws.Range("E41").Value = '  -0.66%  '

$ws.Range("E42").Value = '  -0.88%  '

$ws.Range("B43").Value = 'FraxShare'

$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.56'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -0.39%  '

$ws.Range("B44").Value = 'Quant'

$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.78'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +1.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.804'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +0.53%  '

$ws.Range("D46").Value = '1.765.92'

$ws.Range("E46").Value = '  +0.41%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.27'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0511'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +1.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.45'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +6.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.417'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -0.33%  '

$ws.Range("E51").Value = '  -1.21%  '
